# [improvements] upgraded to Spring Boot 2.1.6
# [json] storeKeys(json,jsonpath,var): extract immediate keys of resolved
#        JSON fragment based on jsonpath.
#
# The "#system" sheet is a hidden lookup table that backs the dropdown
# pickers on the "Scenario" sheet (named ranges point at columns there).
# This edit:
#   1. Inserts a new row for "storeKeys(json,jsonpath,var)" into the
#      alphabetically-sorted "json" column (M), between storeCount and
#      storeValue - shifting only that column's cells down by one.
#   2. Removes the "text" entry from the "target" list (column A),
#      shifting the remaining entries up by one.
#   3. Removes the (now redundant) "text" column (Y) entirely, shifting
#      every column after it (web/webalert/webcookie/ws/ws.async/xml)
#      one letter to the left.
#   4. Updates the named ranges that describe each column's extent to
#      match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. Insert "storeKeys(json,jsonpath,var)" into the json column (M),
#        shifting only column M down starting at row 16 (so that the new
#        row lands alphabetically between storeCount and storeValue).
#        NOTE: Range.Insert()/Delete() on this host shift the *entire*
#        row, so the column-only shift is done by hand, cell by cell.
$jsonCol = 13  # column M
$lastJsonRow = 17
for ($r = $lastJsonRow; $r -ge 16; $r--) {
    $v = $ws.Cells.Item($r, $jsonCol).Value2
    $ws.Cells.Item($r + 1, $jsonCol).Value = $v
}
$ws.Cells.Item(16, $jsonCol).Value = "storeKeys(json,jsonpath,var)"

# --- 2. Remove "text" from the target column (A), shifting A26:A31 up to
#        A25:A30, leaving A31 empty. Again done cell by cell so only
#        column A is touched.
$targetCol = 1  # column A
for ($r = 25; $r -le 30; $r++) {
    $v = $ws.Cells.Item($r + 1, $targetCol).Value2
    $ws.Cells.Item($r, $targetCol).Value = $v
}
$ws.Cells.Item(31, $targetCol).ClearContents()

# --- 3. Delete column Y entirely: shift columns Z..AE left by one
#        (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD), clearing what was
#        column AE. Spans every row that has data in those columns.
$firstShiftCol = 25  # column Y
$lastShiftCol = 30   # column AD (receives old AE)
$clearCol = 31       # column AE (emptied after the shift)
$lastRow = 129
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = $firstShiftCol; $c -le $lastShiftCol; $c++) {
        $v = $ws.Cells.Item($r, $c + 1).Value2
        if ($v -eq $null) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $v
        }
    }
    $ws.Cells.Item($r, $clearCol).ClearContents()
}

# --- 4. Fix up the named ranges that describe each column's extent.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
